$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new rows before existing row 110, shifting old rows 110:122 down to 113:125
$ws.Range("A110:A112").EntireRow.Insert()

# Common (constant) values shared by every row in this data block
$mercadoId   = 6
$mercado     = "Mercado Mayorista Lo Valledor de Santiago"
$region      = "Metropolitana"
$codreg      = 13
$tipo        = "Fruta"
$productoId  = 100104
$producto    = "Frutos de pepita"
$categoriaId = 100104003
$categoria   = "Membrillo"
$variedad    = "Champion"

# ---- New row 110 ----
$ws.Cells.Item(110, 1).Value2  = $mercadoId
$ws.Cells.Item(110, 2).Value2  = $mercado
$ws.Cells.Item(110, 3).Value2  = $region
$ws.Cells.Item(110, 4).Value2  = 44748
$ws.Cells.Item(110, 5).Value2  = $codreg
$ws.Cells.Item(110, 6).Value2  = $tipo
$ws.Cells.Item(110, 7).Value2  = $productoId
$ws.Cells.Item(110, 8).Value2  = $producto
$ws.Cells.Item(110, 9).Value2  = $categoriaId
$ws.Cells.Item(110, 10).Value2 = $categoria
$ws.Cells.Item(110, 11).Value2 = $variedad
$ws.Cells.Item(110, 12).Value2 = "Especial"
$ws.Cells.Item(110, 13).Value2 = 4
$ws.Cells.Item(110, 14).Value2 = 280000
$ws.Cells.Item(110, 15).Value2 = 280000
$ws.Cells.Item(110, 16).Value2 = 280000
$ws.Cells.Item(110, 17).Value2 = "`$/bins (450 kilos)"
$ws.Cells.Item(110, 18).Value2 = "Región de O'Higgins"
$ws.Cells.Item(110, 19).Value2 = 622
$ws.Cells.Item(110, 20).Value2 = 450

# ---- New row 111 ----
$ws.Cells.Item(111, 1).Value2  = $mercadoId
$ws.Cells.Item(111, 2).Value2  = $mercado
$ws.Cells.Item(111, 3).Value2  = $region
$ws.Cells.Item(111, 4).Value2  = 44748
$ws.Cells.Item(111, 5).Value2  = $codreg
$ws.Cells.Item(111, 6).Value2  = $tipo
$ws.Cells.Item(111, 7).Value2  = $productoId
$ws.Cells.Item(111, 8).Value2  = $producto
$ws.Cells.Item(111, 9).Value2  = $categoriaId
$ws.Cells.Item(111, 10).Value2 = $categoria
$ws.Cells.Item(111, 11).Value2 = $variedad
$ws.Cells.Item(111, 12).Value2 = "Primera"
$ws.Cells.Item(111, 13).Value2 = 6
$ws.Cells.Item(111, 14).Value2 = 240000
$ws.Cells.Item(111, 15).Value2 = 240000
$ws.Cells.Item(111, 16).Value2 = 240000
$ws.Cells.Item(111, 17).Value2 = "`$/bins (450 kilos)"
$ws.Cells.Item(111, 18).Value2 = "Región de O'Higgins"
$ws.Cells.Item(111, 19).Value2 = 533
$ws.Cells.Item(111, 20).Value2 = 450

# ---- New row 112 ----
$ws.Cells.Item(112, 1).Value2  = $mercadoId
$ws.Cells.Item(112, 2).Value2  = $mercado
$ws.Cells.Item(112, 3).Value2  = $region
$ws.Cells.Item(112, 4).Value2  = 44748
$ws.Cells.Item(112, 5).Value2  = $codreg
$ws.Cells.Item(112, 6).Value2  = $tipo
$ws.Cells.Item(112, 7).Value2  = $productoId
$ws.Cells.Item(112, 8).Value2  = $producto
$ws.Cells.Item(112, 9).Value2  = $categoriaId
$ws.Cells.Item(112, 10).Value2 = $categoria
$ws.Cells.Item(112, 11).Value2 = $variedad
$ws.Cells.Item(112, 12).Value2 = "Segunda"
$ws.Cells.Item(112, 13).Value2 = 8
$ws.Cells.Item(112, 14).Value2 = 220000
$ws.Cells.Item(112, 15).Value2 = 220000
$ws.Cells.Item(112, 16).Value2 = 220000
$ws.Cells.Item(112, 17).Value2 = "`$/bins (450 kilos)"
$ws.Cells.Item(112, 18).Value2 = "Región de O'Higgins"
$ws.Cells.Item(112, 19).Value2 = 489
$ws.Cells.Item(112, 20).Value2 = 450
